# Update part 3: adjust get_alpha() to monthly return
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003397341679187313
$ws.Range("C2").Value = 0.005060226800726117
$ws.Range("D2").Value = 0.005880069239839422
$ws.Range("E2").Value = 0.005935788986519947
$ws.Range("B3").Value = 1.533100060524297
$ws.Range("C3").Value = 2.730250707311016
$ws.Range("D3").Value = 3.838105272924532
$ws.Range("E3").Value = 4.513695152992129
$ws.Range("B4").Value = -0.009824086724663483
$ws.Range("C4").Value = -0.01215371797881776
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("B5").Value = -2.511590027354631
$ws.Range("C5").Value = -3.179224072623439
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("B6").Value = 0.0003739926898596067
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = -0.007122992673758478
$ws.Range("E6").ClearContents()
$ws.Range("B7").Value = 0.08580150919415963
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = -1.752827438913984
$ws.Range("E7").ClearContents()
$ws.Range("B8").Value = 0.00816532825476327
$ws.Range("C8").Value = 0.007032806108901811
$ws.Range("D8").Value = 0.007338361278754536
$ws.Range("E8").Value = 0.005928761921112338
$ws.Range("B9").Value = 2.885206247699012
$ws.Range("C9").Value = 3.458849286189286
$ws.Range("D9").Value = 4.756818425555082
$ws.Range("E9").Value = 4.871914683516339
$ws.Range("B10").Value = -0.008855798268871057
$ws.Range("C10").Value = -0.01172040886061532
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("B11").Value = -2.246623680712763
$ws.Range("C11").Value = -3.068733526964793
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("B12").Value = -0.006952026373694971
$ws.Range("C12").Value = -0.01068788666429326
$ws.Range("D12").Value = -0.008176746263975057
$ws.Range("E12").ClearContents()
$ws.Range("B13").Value = -1.310774187180516
$ws.Range("C13").Value = -1.975209557545742
$ws.Range("D13").Value = -2.109753324264922
$ws.Range("E13").ClearContents()
$ws.Range("B14").Value = 0.006401444703493888
$ws.Range("C14").Value = 0.007379937413733258
$ws.Range("D14").Value = 0.007351262216361159
$ws.Range("E14").Value = 0.005638523698822219
$ws.Range("B15").Value = 2.071757768678078
$ws.Range("C15").Value = 3.606705359259601
$ws.Range("D15").Value = 4.840882952723406
$ws.Range("E15").Value = 4.786377614314292
$ws.Range("B16").Value = -0.009212670453859604
$ws.Range("C16").Value = -0.01236815867227868
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("B17").Value = -2.305547223198694
$ws.Range("C17").Value = -3.112649879175956
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("B18").Value = -0.005321318143397036
$ws.Range("C18").Value = -0.01023302863420526
$ws.Range("D18").Value = -0.007481742676284792
$ws.Range("E18").Value = -0.00235109639791297
$ws.Range("B19").Value = -0.9861119719109196
$ws.Range("C19").Value = -1.956883898333945
$ws.Range("D19").Value = -1.967730362415597
$ws.Range("E19").Value = -0.7436538190882448
$ws.Range("B20").Value = 0.00659235594305914
$ws.Range("C20").Value = 0.008101624932892503
$ws.Range("D20").Value = 0.007750688364177777
$ws.Range("E20").Value = 0.006350500135489026
$ws.Range("B21").Value = 2.138047777371845
$ws.Range("C21").Value = 3.834875614454537
$ws.Range("D21").Value = 4.855638079197258
$ws.Range("E21").Value = 5.032686349931778
$ws.Range("B22").Value = -0.01014972583053267
$ws.Range("C22").Value = -0.01333983249292729
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("B23").Value = -2.52614267841462
$ws.Range("C23").Value = -3.325316466907414
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("B24").Value = -0.004493840652280583
$ws.Range("C24").Value = -0.01031677158639364
$ws.Range("D24").Value = -0.007749540885006581
$ws.Range("E24").Value = -0.003238046131010818
$ws.Range("B25").Value = -0.8272399589667041
$ws.Range("C25").Value = -1.929615649771915
$ws.Range("D25").Value = -1.958322078928385
$ws.Range("E25").Value = -1.037005651016436

Write-Host "done"
